$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.27738566666667
$ws.Range("H2").Value = 30.832157
$ws.Range("I2").Value = 0.3571200664977529
$ws.Range("J2").Value = 0.3571200664977529
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.435593666666667
$ws.Range("N2").Value = 13.306781
$ws.Range("O2").Value = 0.2725293883759869
$ws.Range("P2").Value = 0.2725293883759869
$ws.Range("Q2").Value = 45.58630677295744
$ws.Range("R2").Value = 410.276760956617
$ws.Range("S2").Value = 0.09732571329942435
$ws.Range("T2").Value = 0.09732571329942437

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.27738566666667
$ws.Range("H3").Value = 30.832157
$ws.Range("I3").Value = 0.3571200664977529
$ws.Range("J3").Value = 0.3571200664977529
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 5.974929333333333
$ws.Range("N3").Value = 17.924788
$ws.Range("O3").Value = 0.3671084321902667
$ws.Range("P3").Value = 0.3671084321902667
$ws.Range("Q3").Value = 61.40665308974621
$ws.Range("R3").Value = 552.6598778077159
$ws.Range("S3").Value = 0.1311017877156738
$ws.Range("T3").Value = 0.1311017877156739

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.27738566666667
$ws.Range("H4").Value = 30.832157
$ws.Range("I4").Value = 0.3571200664977529
$ws.Range("J4").Value = 0.3571200664977529
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.398978
$ws.Range("N4").Value = 10.196934
$ws.Range("O4").Value = 0.2088381995863842
$ws.Range("P4").Value = 0.2088381995863842
$ws.Range("Q4").Value = 34.93260777851533
$ws.Range("R4").Value = 314.393470006638
$ws.Range("S4").Value = 0.07458031172356053
$ws.Range("T4").Value = 0.07458031172356054

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.27738566666667
$ws.Range("H5").Value = 30.832157
$ws.Range("I5").Value = 0.3571200664977529
$ws.Range("J5").Value = 0.3571200664977529
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.466151666666667
$ws.Range("N5").Value = 7.398455
$ws.Range("O5").Value = 0.1515239798473622
$ws.Range("P5").Value = 0.1515239798473622
$ws.Range("Q5").Value = 25.34559179082611
$ws.Range("R5").Value = 228.110326117435
$ws.Range("S5").Value = 0.05411225375909415
$ws.Range("T5").Value = 0.05411225375909416

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.987132666666667
$ws.Range("H6").Value = 5.961398
$ws.Range("I6").Value = 0.06904917000064482
$ws.Range("J6").Value = 0.06904917000064482
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 4.435593666666667
$ws.Range("N6").Value = 13.306781
$ws.Range("O6").Value = 0.2725293883759869
$ws.Range("P6").Value = 0.2725293883759869
$ws.Range("Q6").Value = 8.814113071093111
$ws.Range("R6").Value = 79.327017639838
$ws.Range("S6").Value = 0.01881792806814528
$ws.Range("T6").Value = 0.01881792806814528

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.987132666666667
$ws.Range("H7").Value = 5.961398
$ws.Range("I7").Value = 0.06904917000064482
$ws.Range("J7").Value = 0.06904917000064482
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 5.974929333333333
$ws.Range("N7").Value = 17.924788
$ws.Range("O7").Value = 0.3671084321902667
$ws.Range("P7").Value = 0.3671084321902667
$ws.Range("Q7").Value = 11.87297725929155
$ws.Range("R7").Value = 106.856795333624
$ws.Range("S7").Value = 0.02534853254297592
$ws.Range("T7").Value = 0.02534853254297592

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.987132666666667
$ws.Range("H8").Value = 5.961398
$ws.Range("I8").Value = 0.06904917000064482
$ws.Range("J8").Value = 0.06904917000064482
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.398978
$ws.Range("N8").Value = 10.196934
$ws.Range("O8").Value = 0.2088381995863842
$ws.Range("P8").Value = 0.2088381995863842
$ws.Range("Q8").Value = 6.754220217081333
$ws.Range("R8").Value = 60.787981953732
$ws.Range("S8").Value = 0.01442010434586884
$ws.Range("T8").Value = 0.01442010434586884

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.987132666666667
$ws.Range("H9").Value = 5.961398
$ws.Range("I9").Value = 0.06904917000064482
$ws.Range("J9").Value = 0.06904917000064482
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.466151666666667
$ws.Range("N9").Value = 7.398455
$ws.Range("O9").Value = 0.1515239798473622
$ws.Range("P9").Value = 0.1515239798473622
$ws.Range("Q9").Value = 4.900570537787778
$ws.Range("R9").Value = 44.10513484009
$ws.Range("S9").Value = 0.01046260504365479
$ws.Range("T9").Value = 0.01046260504365479

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 15.70994033333333
$ws.Range("H10").Value = 47.129821
$ws.Range("I10").Value = 0.5458912527445677
$ws.Range("J10").Value = 0.5458912527445677
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.435593666666667
$ws.Range("N10").Value = 13.306781
$ws.Range("O10").Value = 0.2725293883759869
$ws.Range("P10").Value = 0.2725293883759869
$ws.Range("Q10").Value = 69.68291184624455
$ws.Range("R10").Value = 627.1462066162011
$ws.Range("S10").Value = 0.1487714092302783
$ws.Range("T10").Value = 0.1487714092302783

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 15.70994033333333
$ws.Range("H11").Value = 47.129821
$ws.Range("I11").Value = 0.5458912527445677
$ws.Range("J11").Value = 0.5458912527445677
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 5.974929333333333
$ws.Range("N11").Value = 17.924788
$ws.Range("O11").Value = 0.3671084321902667
$ws.Range("P11").Value = 0.3671084321902667
$ws.Range("Q11").Value = 93.86578332254977
$ws.Range("R11").Value = 844.792049902948
$ws.Range("S11").Value = 0.2004012819414389
$ws.Range("T11").Value = 0.2004012819414389

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 15.70994033333333
$ws.Range("H12").Value = 47.129821
$ws.Range("I12").Value = 0.5458912527445677
$ws.Range("J12").Value = 0.5458912527445677
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.398978
$ws.Range("N12").Value = 10.196934
$ws.Range("O12").Value = 0.2088381995863842
$ws.Range("P12").Value = 0.2088381995863842
$ws.Range("Q12").Value = 53.39774157431267
$ws.Range("R12").Value = 480.5796741688141
$ws.Range("S12").Value = 0.1140029463931314
$ws.Range("T12").Value = 0.1140029463931314

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 15.70994033333333
$ws.Range("H13").Value = 47.129821
$ws.Range("I13").Value = 0.5458912527445677
$ws.Range("J13").Value = 0.5458912527445677
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.466151666666667
$ws.Range("N13").Value = 7.398455
$ws.Range("O13").Value = 0.1515239798473622
$ws.Range("P13").Value = 0.1515239798473622
$ws.Range("Q13").Value = 38.74309553628389
$ws.Range("R13").Value = 348.687859826555
$ws.Range("S13").Value = 0.08271561517971919
$ws.Range("T13").Value = 0.08271561517971919

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.8040576666666667
$ws.Range("H14").Value = 2.412173
$ws.Range("I14").Value = 0.02793951075703474
$ws.Range("J14").Value = 0.02793951075703475
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4.435593666666667
$ws.Range("N14").Value = 13.306781
$ws.Range("O14").Value = 0.2725293883759869
$ws.Range("P14").Value = 0.2725293883759869
$ws.Range("Q14").Value = 3.566473093901445
$ws.Range("R14").Value = 32.098257845113
$ws.Range("S14").Value = 0.007614337778138985
$ws.Range("T14").Value = 0.007614337778138986

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.8040576666666667
$ws.Range("H15").Value = 2.412173
$ws.Range("I15").Value = 0.02793951075703474
$ws.Range("J15").Value = 0.02793951075703475
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 5.974929333333333
$ws.Range("N15").Value = 17.924788
$ws.Range("O15").Value = 0.3671084321902667
$ws.Range("P15").Value = 0.3671084321902667
$ws.Range("Q15").Value = 4.804187738258221
$ws.Range("R15").Value = 43.237689644324
$ws.Range("S15").Value = 0.01025682999017812
$ws.Range("T15").Value = 0.01025682999017812

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.8040576666666667
$ws.Range("H16").Value = 2.412173
$ws.Range("I16").Value = 0.02793951075703474
$ws.Range("J16").Value = 0.02793951075703475
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 3.398978
$ws.Range("N16").Value = 10.196934
$ws.Range("O16").Value = 0.2088381995863842
$ws.Range("P16").Value = 0.2088381995863842
$ws.Range("Q16").Value = 2.732974319731333
$ws.Range("R16").Value = 24.596768877582
$ws.Range("S16").Value = 0.005834837123823552
$ws.Range("T16").Value = 0.005834837123823553

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.8040576666666667
$ws.Range("H17").Value = 2.412173
$ws.Range("I17").Value = 0.02793951075703474
$ws.Range("J17").Value = 0.02793951075703475
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.466151666666667
$ws.Range("N17").Value = 7.398455
$ws.Range("O17").Value = 0.1515239798473622
$ws.Range("P17").Value = 0.1515239798473622
$ws.Range("Q17").Value = 1.982928154746111
$ws.Range("R17").Value = 17.846353392715
$ws.Range("S17").Value = 0.004233505864894092
$ws.Range("T17").Value = 0.004233505864894093
